# Generate Report for Handoff
# Updates the localization-status report: the b.md file has now been
# handed off for localization (zh-cn / de-de), so its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff" across the
# Overview sheet and each language sheet, along with the new handoff
# file name / timestamp and (for zh-cn/de-de) a stale-handback warning.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: row for b.md (row 3) ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-09-02 12:43:37"

# ---- zh-cn sheet: row for b.md (row 3) ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "False"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws.Range("H3").Value = "2016-09-02 12:43:32"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/75f3cc01b56debd7e68f648b2fe85964512a4d91/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f824cdc2512469dbeb9c158f9b0b98b973714e04/e2e/b.md."
$ws.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet: row for b.md (row 3) ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "False"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws.Range("H3").Value = "2016-09-02 12:43:37"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/75f3cc01b56debd7e68f648b2fe85964512a4d91/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f824cdc2512469dbeb9c158f9b0b98b973714e04/e2e/b.md."
$ws.Columns.Item(16).ColumnWidth = 39.166666666666664
